# Reorder the staff names listed in column A (rows 2-35) of the "staff" sheet.
# The list of names itself is unchanged; only the row each name occupies changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("staff")

$names = @(
    'Harris, Jack',
    'Allan, Steven',
    'Harris, Lily',
    'Hardacre, Trevor',
    'Carr, Isaac',
    'Burgess, Isaac',
    'Rutherford, Alexander',
    'Hamilton, Lisa',
    'Cornish, Nathan',
    'Terry, Sue',
    'Randall, Jason',
    'Jackson, Max',
    'Baker, Kevin',
    'Smith, Connor',
    'Howard, Faith',
    'Peters, Megan',
    'Smith, Jonathan',
    'BYRNE, JOHN (SHANE)',
    'Miller, Alison',
    'Mackay, Claire',
    'Glover, Michael',
    'Marshall, Vanessa',
    'Jones, James',
    'Clarkson, Peter',
    'Edmunds, Colin',
    'Carr, Connor',
    'Wilkins, Wendy',
    'Paige, Jennifer',
    'Powell, Blake',
    'Edmunds, Madeleine',
    'Nolan, Stewart',
    'Duncan, Michelle',
    'Newman, Andrew',
    'MacLeod, Julia'
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}
